function Format-PyFloat($v) {
    # Produce the shortest decimal string that round-trips to the exact
    # same double value (mirrors Python's repr()/str() behavior for floats),
    # including a trailing ".0" for integral values. Avoids .NET's default
    # double->string conversion, which in this runtime does not reliably
    # round-trip (and can also emit scientific notation).
    for ($p = 1; $p -le 17; $p++) {
        $fmt = "G" + $p
        $s = $v.ToString($fmt)
        $parsed = [double]$s
        if ($parsed -eq $v) {
            if ($s -match '[eE]') {
                # Re-render in fixed-point notation (Python only switches to
                # exponent form far outside the magnitudes used here).
                for ($q = 0; $q -le 17; $q++) {
                    $cand = $v.ToString("F$q")
                    if ([double]$cand -eq $v) {
                        $s = $cand
                        break
                    }
                }
            }
            if ($s -notmatch '[.eE]') {
                $s = "{0}.0" -f $s
            }
            return $s
        }
    }
    return $v.ToString("G17")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Left (L) block: add "Running Total" header + running totals in column B ---
$ws.Range("B3").Value = "Running Total"

$leftValues = @()
$leftRunning = 0
for ($r = 4; $r -le 38; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $leftRunning = $leftRunning + $a
    $ws.Cells.Item($r, 2).Value = $leftRunning
    $leftValues += $leftRunning
}

# --- Right (R) block: add "Running Total" header + running totals in column B ---
$ws.Range("B39").Value = "Running Total"

$rightValues = @()
$rightRunning = 0
for ($r = 40; $r -le 46; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $rightRunning = $rightRunning + $a
    $ws.Cells.Item($r, 2).Value = $rightRunning
    $rightValues += $rightRunning
}

# --- Relabel the raw-data rows ---
$ws.Range("A47").Value = "L ~ Raw"
$ws.Range("A49").Value = "R ~ Raw"

# --- Append running-total summary rows (51-54) ---
$leftParts = @()
foreach ($v in $leftValues) {
    $leftParts += Format-PyFloat $v
}
$leftSummary = [string]::Join(", ", $leftParts)

$rightParts = @()
foreach ($v in $rightValues) {
    $rightParts += Format-PyFloat $v
}
$rightSummary = [string]::Join(", ", $rightParts)

$ws.Range("A51").Value = "L ~ Running Totals"
$ws.Range("A52").Value = $leftSummary
$ws.Range("A53").Value = "R ~ Running Totals"
$ws.Range("A54").Value = $rightSummary
